# Automatische test-sync: 2025-08-04 20:40:50
# Adds a new test-mail row (row 19) to the "Logs" sheet, extends the
# conditional-formatting ranges that cover the data rows, and bumps the
# "Inkoop / Bestellingen" tally on the "Dashboard" sheet from 4 to 5.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new row of data (row 19) -------------------------------
$logs.Range("A19").Value = "Is dit artikel nog op voorraad?"
$logs.Range("B19").Value = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value = "Testmail #7: Is dit artikel nog op voorraad?"
$logs.Range("D19").Value = "Inkoop / Bestellingen"
$logs.Range("E19").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F19").Value = "2025-08-04 20:40:33"
$logs.Range("G19").Value = "Ja"
$logs.Range("H19").Value = "Ja"
$logs.Range("I19").Value = "Nee"
$logs.Range("J19").Value = "Nee"

# --- 2. Extend the conditional formatting ranges to include row 19 --------
# Each column's conditional-formatting block (D, G, H, I, J) previously
# covered rows 2-18; it now needs to cover rows 2-19. Modifying the
# AppliesTo range of any one rule in a block updates the whole block.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "18")
    $newRange = $logs.Range($col + "2:" + $col + "19")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the Dashboard summary count for "Inkoop / Bestellingen" ----
$dashboard.Range("B3").Value = 5
